$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.367.03'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.846.26'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9982'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.09'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6275'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2900'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.48'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07738'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.846.46'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.21%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.992'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6803'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.05'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.193'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.385.61'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.97'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.34'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9992'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.53'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.436'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1368'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.52'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06500'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +16.06%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.408'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.481'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.081'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.831'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6976'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.581'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.264.82'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.19%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.832'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.19%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01835'
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.751'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.24%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9179'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.28%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9990'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.008.60'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.43'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Aptos'
$ws.Range("B46").Style = "Normal"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C46").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.083'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.49%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C47").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.723'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("B48").Style = "Normal"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("C48").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000117'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1163'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("B50").Style = "Normal"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C50").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3959'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("B51").Style = "Normal"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C51").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.966'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.50%  '
$ws.Range("E51").Style = "Normal"
